# Add team record columns (Wins, Losses, Ties) to the BAL_2018 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy formatting from the existing last header cell (AC1)
# so the new headers match the bold/border/centered style used by the other headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2 through 58): every player on the roster shares the same
# team win/loss/tie record for the season.
$lastRow = 58
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 47    # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 115   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0     # AF -> Ties
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-$lastRow"
